$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 2.25
$ws.Range("X3").Value = 9.5
$ws.Range("Y3").Value = 8.5
$ws.Range("AB3").Value = 23
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 13
$ws.Range("AO3").Value = 10
$ws.Range("AP3").Value = 19
$ws.Range("AX3").Value = 21

# Row 4 updates
$ws.Range("G4").Value = 2.52
$ws.Range("I4").Value = 2.35
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 3.55
$ws.Range("Q4").Value = 1.75
$ws.Range("U4").Value = 1.61
$ws.Range("V4").Value = 2.25
$ws.Range("W4").Value = 8
$ws.Range("AA4").Value = 16.5
$ws.Range("AB4").Value = 22
$ws.Range("AC4").Value = 11.5
$ws.Range("AH4").Value = 7.7
$ws.Range("AJ4").Value = 8
$ws.Range("AK4").Value = 19.5
